# "Handling of PCF valves"
#
# Changes applied here (per the target diff):
#   1. Shared string "GenValve: Std" -> "GenericPCFValve: Std"
#      (used by cells L2 and L3 on the "All pipelines" sheet).
#   2. Scroll the view so column B leads (topLeftCell="B1") and move the
#      selection from L4 to L9 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the valve description used in column L (rows 2 and 3 both
#    point at the same shared string, so updating one cell's value
#    rewrites the shared string in place for both).
$ws.Range("L2").Value = "GenericPCFValve: Std"
$ws.Range("L3").Value = "GenericPCFValve: Std"

# 2/3. Scroll so column B leads the view, then move the selection to L9.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("L9").Select()
